# Update the "caseId" row (row 2) of the first table (RS-RR schema table):
#  - Column 5 ("Description"): rewrite the explanatory text.
#  - Column 6 ("Exemple"): rewrite the example value.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Description cell -------------------------------------------------
$descCell = $t.Cell(2, 5)
$br = [char]11   # manual line break (renders as <w:br/> when written back)

$newDescription = (
    "Identifiant partagé de l'affaire/dossier, généré une seule fois par le système du partenaire qui recoit la primo-demande de secours (créateur du dossier). " + $br +
    "Il est valorisé comme suit lors de sa création : " + $br +
    "{pays}.{domaine}.{organisation}.{senderCaseId}" + $br + $br +
    "Il doit pouvoir être généré de façon décentralisée et ne présenter aucune ambiguïté." + $br +
    " Il doit être unique dans l'ensemble des systèmes : le numéro de dossier fourni par celui qui génère l'identifiant partagé doit donc être un numéro unique dans son système."
)

$descCell.Range.Text = $newDescription

# --- Exemple cell -------------------------------------------------------
$exampleCell = $t.Cell(2, 6)
$exampleCell.Range.Text = "fr.health.samu440.DRFR15440241550012"
